$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.990.58"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "2.908.16"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'568.36"
$ws.Range("E5").Value = "  -3.35%  "
$ws.Range("D6").Value = "'143.75"
$ws.Range("E6").Value = "  -1.67%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "2.905.27"
$ws.Range("E8").Value = "  -0.32%  "
$ws.Range("D9").Value = "'0.499"
$ws.Range("E9").Value = "  -1.48%  "
$ws.Range("D10").Value = "'6.94"
$ws.Range("E10").Value = "  -0.81%  "
$ws.Range("E12").Value = "  -1.60%  "
$ws.Range("E13").Value = "  -0.43%  "
$ws.Range("D14").Value = "'32.51"
$ws.Range("E14").Value = "  +0.12%  "
$ws.Range("E15").Value = "  +0.22%  "
$ws.Range("D16").Value = "3.389.49"
$ws.Range("E16").Value = "  -0.37%  "
$ws.Range("D17").Value = "61.950.48"
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").Value = "2.906.35"
$ws.Range("E18").Value = "  -0.49%  "
$ws.Range("E19").Value = "  -1.65%  "
$ws.Range("D20").Value = "'429.79"
$ws.Range("E20").Value = "  -1.01%  "
$ws.Range("E21").Value = "  -3.02%  "
$ws.Range("D22").Value = "'0.652"
$ws.Range("E22").Value = "  -1.33%  "
$ws.Range("D23").Value = "'6.87"
$ws.Range("D24").Value = "'78.82"
$ws.Range("E24").Value = "  -2.64%  "
$ws.Range("D25").Value = "'12.00"
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("D26").Value = "'10.28"
$ws.Range("E26").Value = "  -6.33%  "
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("E28").Value = "  -3.46%  "
$ws.Range("E29").Value = "  +9.46%  "
$ws.Range("E31").Value = "  -2.70%  "
$ws.Range("E32").Value = "  -5.11%  "
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("E34").Value = "  -3.02%  "
$ws.Range("D35").Value = "'25.63"
$ws.Range("E35").Value = "  -1.63%  "
$ws.Range("D36").Value = "'0.953"
$ws.Range("E36").Value = "  -2.46%  "
$ws.Range("E37").Value = "  -2.55%  "
$ws.Range("D38").Value = "'48.83"
$ws.Range("E38").Value = "  -0.93%  "
$ws.Range("D39").Value = "'2.90"
$ws.Range("E39").Value = "  -5.86%  "
$ws.Range("E40").Value = "  -4.86%  "
$ws.Range("E41").Value = "  -0.87%  "
$ws.Range("E42").Value = "  +5.55%  "
$ws.Range("E43").Value = "  -2.52%  "
$ws.Range("D44").Value = "'0.267"
$ws.Range("E44").Value = "  -2.77%  "
$ws.Range("D45").Value = "2.713.44"
$ws.Range("E45").Value = "  +0.51%  "
$ws.Range("D46").Value = "'132.94"
$ws.Range("E46").Value = "  -1.52%  "
$ws.Range("E47").Value = "  -0.64%  "
$ws.Range("D48").Value = "'349.69"
$ws.Range("E48").Value = "  +0.84%  "
$ws.Range("E50").Value = "  -1.29%  "
$ws.Range("D51").Value = "'0.000210"
$ws.Range("E51").Value = "  +10.34%  "
